$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Re-sort / correct the "Periodo Mora" rows for worker 74183766 (JAVIER ALONSO MARIN PATIÑO) ---
# Before: row16=2303/79200, row17=2302/108000, row18=2301/46800, row19(other worker)=2503/9490
# After:  row16=2301/46800, row17=2302/108000, row18=2303/79200, row19 removed entirely
#
# First, promote row 18's formatting to the "closing" (bottom-border) style that row 19 currently has,
# since after row 19 is deleted, row 18 becomes the last row of this worker's block.
$ws.Range("B19:J19").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Now fix the period / value data so the three remaining rows read 2301, 2302, 2303 in order.
$ws.Range("E16").Value = "2301"
$ws.Range("F16").Value = 46800

$ws.Range("E17").Value = "2302"
$ws.Range("F17").Value = 108000

$ws.Range("E18").Value = "2303"
$ws.Range("F18").Value = 79200
$ws.Range("G18").Value = 2700000

# --- 2. Remove the second worker's record (73150830 / SHARIFF ANTONIO CABARCAS RAMOS) entirely ---
$ws.Rows.Item(19).Delete()

# --- 3. Update the summary counters / totals at the top of the sheet ---
$ws.Range("C13").Value = 1        # Cant. Trabajadores: 2 -> 1
$ws.Range("F13").Value = 3        # Cant. Periodos: 4 -> 3
$ws.Range("E11").Value = 234000   # VALOR MORA total: 243490 -> 234000

# --- 4. Column D ("Nombre Trabajador") no longer needs to fit the long removed name, re-fit it ---
$ws.Columns.Item(4).AutoFit()
